# "Tomada de testes" — replace the real customer file paths / passwords
# with generic test ("tomada de testes") placeholder data, tighten the
# RIGHT() formulas to match the new (shorter) filename length, re-apply
# column autofit-like widths, add an AutoFilter over the header row (plus
# its companion hidden _FilterDatabase defined name), and move the active
# selection / window chrome the way the author's Excel session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: CAMINHO ARQUIVO (full path) -------------------------------
$ws.Range("A2").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\arquivo01.pdf"
$ws.Range("A3").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\arquivo02.pdf"
$ws.Range("A4").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\arquivo03.pdf"
$ws.Range("A5").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\arquivo04.pdf"

# --- Column B: SENHA ARQUIVO (password) -----------------------------------
$ws.Range("B2").Value = 123
$ws.Range("B3").Value = 123
$ws.Range("B4").Value = 123
$ws.Range("B5").Value = 123

# --- Column C: NOME ARQUIVO (derived filename, RIGHT(path,16) -> RIGHT(path,13)) ---
$ws.Range("C2").Formula = "=RIGHT(A2,13)"
$ws.Range("C3:C5").Formula = "=RIGHT(A3,13)"

# B3:B5 picked up C3's border formatting (thin box, no top border) instead
# of the original full-box border used by A3:A5/B2.
$ws.Range("C3").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)  # xlPasteFormats

# --- Column widths (re-fit to the new, shorter content) -------------------
$ws.Columns.Item(1).ColumnWidth = 60.3
$ws.Columns.Item(2).ColumnWidth = 19.8
$ws.Columns.Item(3).ColumnWidth = 19.6

# --- AutoFilter over the header row + its hidden defined name -------------
$null = $ws.Range("A1:C1").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:C1"))
$filterName.Visible = $false

# --- Selection / window chrome --------------------------------------------
$null = $ws.Range("F15").Select()

$win = $wb.Windows.Item(1)
$win.Left = -105
$win.Top = -105
$win.Width = 23250
$win.Height = 12450
